$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H38").Value = 1935
$ws_ALC.Range("I38").Value = 927
$ws_ALC.Range("J38").Value = 9999
$ws_ALC.Range("K38").Value = 2781
$ws_ALC.Range("L38").Value = 29997
$ws_ALC.Range("M38").Value = -2409
$ws_ALC.Range("N38").Value = -30741
$ws_ALC.Range("H40").Value = 2192.3076
$ws_ALC.Range("I40").Value = 1800
$ws_ALC.Range("J40").Value = 2437.5
$ws_ALC.Range("K40").Value = 1800
$ws_ALC.Range("L40").Value = 2437.5
$ws_ALC.Range("M40").Value = -1625
$ws_ALC.Range("N40").Value = -2787.5
$ws_ALC.Range("H58").Value = 738.3333
$ws_ALC.Range("I58").Value = 738.3333
$ws_ALC.Range("J58").Value = 0
$ws_ALC.Range("K58").Value = 2214.9999
$ws_ALC.Range("L58").Value = 0
$ws_ALC.Range("M58").ClearContents()
$ws_ALC.Range("N58").Value = -2064.9999
$ws_ALC.Range("H62").Value = 4811.75
$ws_ALC.Range("I62").Value = 4811.75
$ws_ALC.Range("J62").Value = 0
$ws_ALC.Range("K62").Value = 4811.75
$ws_ALC.Range("L62").Value = 0
$ws_ALC.Range("M62").ClearContents()
$ws_ALC.Range("N62").Value = -4187.75
$ws_ALC.Range("H65").Value = 4811.75
$ws_ALC.Range("I65").Value = 4811.75
$ws_ALC.Range("J65").Value = 0
$ws_ALC.Range("K65").Value = 24058.75
$ws_ALC.Range("L65").Value = 0
$ws_ALC.Range("M65").Value = -20938.75
$ws_ALC.Range("H70").Value = 4700.1177
$ws_ALC.Range("J70").Value = 5026.8
$ws_ALC.Range("L70").Value = 15080.4
$ws_ALC.Range("N70").Value = -15620.4
$ws_ALC.Range("H73").Value = 4700.1177
$ws_ALC.Range("J73").Value = 5026.8
$ws_ALC.Range("L73").Value = 15080.4
$ws_ALC.Range("N73").Value = -16952.4
$ws_ALC.Range("H132").Value = 1434.9
$ws_ALC.Range("I132").Value = 1480.76
$ws_ALC.Range("J132").Value = 1205.6
$ws_ALC.Range("K132").Value = 4442.28
$ws_ALC.Range("L132").Value = 3616.8
$ws_ALC.Range("M132").Value = -1912.28
$ws_ALC.Range("N132").Value = -8676.799999999999
$ws_ALC.Range("H137").Value = 2138.5386
$ws_ALC.Range("I137").Value = 2222.4443
$ws_ALC.Range("K137").Value = 6667.3329
$ws_ALC.Range("M137").Value = -4117.3329

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H45").Value = 1974.5
$ws_ARM.Range("I45").Value = 1974.5
$ws_ARM.Range("J45").Value = 0
$ws_ARM.Range("K45").Value = 1974.5
$ws_ARM.Range("L45").Value = 0
$ws_ARM.Range("M45").ClearContents()
$ws_ARM.Range("N45").Value = -1597.5
$ws_ARM.Range("H61").Value = 3144.3333
$ws_ARM.Range("J61").Value = 0
$ws_ARM.Range("L61").Value = 0
$ws_ARM.Range("N61").ClearContents()
$ws_ARM.Range("H102").Value = 1194.8
$ws_ARM.Range("I102").Value = 993.5
$ws_ARM.Range("J102").Value = 2000
$ws_ARM.Range("K102").Value = 993.5
$ws_ARM.Range("L102").Value = 2000
$ws_ARM.Range("M102").Value = 628.5
$ws_ARM.Range("N102").Value = -5244
$ws_ARM.Range("H122").Value = 437320.8
$ws_ARM.Range("I122").Value = 771494.1
$ws_ARM.Range("K122").Value = 2314482.3
$ws_ARM.Range("M122").Value = -2312032.3
$ws_ARM.Range("H132").Value = 1067.7333
$ws_ARM.Range("I132").Value = 1113.1666
$ws_ARM.Range("J132").Value = 886
$ws_ARM.Range("K132").Value = 3339.4998
$ws_ARM.Range("L132").Value = 2658
$ws_ARM.Range("M132").Value = -809.4998000000001
$ws_ARM.Range("N132").Value = -7718
$ws_ARM.Range("H136").Value = 3144.3333
$ws_ARM.Range("J136").Value = 0
$ws_ARM.Range("L136").Value = 0
$ws_ARM.Range("N136").ClearContents()

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 2630.575
$ws_CRP.Range("J31").Value = 2442.182
$ws_CRP.Range("L31").Value = 2442.182
$ws_CRP.Range("N31").Value = -3032.182
$ws_CRP.Range("H34").Value = 2630.575
$ws_CRP.Range("J34").Value = 2442.182
$ws_CRP.Range("L34").Value = 2442.182
$ws_CRP.Range("N34").Value = -2846.182
$ws_CRP.Range("H134").Value = 1732.3773
$ws_CRP.Range("I134").Value = 1483.55
$ws_CRP.Range("K134").Value = 4450.65
$ws_CRP.Range("M134").Value = -1915.65

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H38").Value = 115.4
$ws_CUL.Range("I38").Value = 64.833336
$ws_CUL.Range("J38").Value = 191.25
$ws_CUL.Range("K38").Value = 194.500008
$ws_CUL.Range("L38").Value = 573.75
$ws_CUL.Range("M38").Value = 152.499992
$ws_CUL.Range("N38").Value = -1267.75
$ws_CUL.Range("H132").Value = 8982.4375
$ws_CUL.Range("I132").Value = 9562.143
$ws_CUL.Range("K132").Value = 86059.287
$ws_CUL.Range("M132").Value = -83529.287

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H80").Value = 2984.8
$ws_GSM.Range("I80").Value = 2474.6667
$ws_GSM.Range("K80").Value = 2474.6667
$ws_GSM.Range("M80").Value = -1476.6667
$ws_GSM.Range("H83").Value = 2984.8
$ws_GSM.Range("I83").Value = 2474.6667
$ws_GSM.Range("K83").Value = 12373.3335
$ws_GSM.Range("M83").Value = -7381.333500000001
$ws_GSM.Range("H102").Value = 1411.7391
$ws_GSM.Range("I102").Value = 598.3333
$ws_GSM.Range("K102").Value = 598.3333
$ws_GSM.Range("M102").Value = 1023.6667
$ws_GSM.Range("H113").Value = 2998.5
$ws_GSM.Range("J113").Value = 2997
$ws_GSM.Range("L113").Value = 2997
$ws_GSM.Range("N113").Value = -7337
$ws_GSM.Range("H122").Value = 35398.1
$ws_GSM.Range("I122").Value = 1948.6957
$ws_GSM.Range("K122").Value = 5846.0871
$ws_GSM.Range("M122").Value = -3396.0871
$ws_GSM.Range("H123").Value = 24250.572
$ws_GSM.Range("J123").Value = 24250.572
$ws_GSM.Range("L123").Value = 24250.572
$ws_GSM.Range("N123").Value = -29150.572
$ws_GSM.Range("H132").Value = 2900.0833
$ws_GSM.Range("I132").Value = 2864.6667
$ws_GSM.Range("K132").Value = 8594.000100000001
$ws_GSM.Range("M132").Value = -6064.000100000001

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H16").Value = 18190.637
$ws_LTW.Range("I16").Value = 19332.666
$ws_LTW.Range("J16").Value = 17762.375
$ws_LTW.Range("K16").Value = 19332.666
$ws_LTW.Range("L16").Value = 17762.375
$ws_LTW.Range("M16").Value = -19162.666
$ws_LTW.Range("N16").Value = -18102.375
$ws_LTW.Range("H82").Value = 1662
$ws_LTW.Range("I82").Value = 1303.3334
$ws_LTW.Range("K82").Value = 1303.3334
$ws_LTW.Range("M82").Value = -942.3334
$ws_LTW.Range("H85").Value = 1662
$ws_LTW.Range("I85").Value = 1303.3334
$ws_LTW.Range("K85").Value = 1303.3334
$ws_LTW.Range("M85").Value = -55.33339999999998
$ws_LTW.Range("H133").Value = 12500
$ws_LTW.Range("J133").Value = 12500
$ws_LTW.Range("L133").Value = 12500
$ws_LTW.Range("N133").Value = -17560

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H74").Value = 60898.668
$ws_WVR.Range("J74").Value = 60898.668
$ws_WVR.Range("L74").Value = 60898.668
$ws_WVR.Range("N74").Value = -62770.668
$ws_WVR.Range("H77").Value = 60898.668
$ws_WVR.Range("J77").Value = 60898.668
$ws_WVR.Range("L77").Value = 182696.004
$ws_WVR.Range("N77").Value = -192056.004
$ws_WVR.Range("H81").Value = 8792.786
$ws_WVR.Range("J81").Value = 9249.950000000001
$ws_WVR.Range("L81").Value = 18499.9
$ws_WVR.Range("N81").Value = -20621.9
$ws_WVR.Range("H84").Value = 8792.786
$ws_WVR.Range("J84").Value = 9249.950000000001
$ws_WVR.Range("L84").Value = 92499.5
$ws_WVR.Range("N84").Value = -103107.5
$ws_WVR.Range("H110").Value = 0
$ws_WVR.Range("J110").Value = 0
$ws_WVR.Range("L110").ClearContents()
$ws_WVR.Range("N110").Value = 0
$ws_WVR.Range("H122").Value = 1105.75
$ws_WVR.Range("I122").Value = 1067.2858
$ws_WVR.Range("K122").Value = 3201.8574
$ws_WVR.Range("M122").Value = -751.8574000000003
$ws_WVR.Range("H132").Value = 29636.545
$ws_WVR.Range("I132").Value = 37175.54
$ws_WVR.Range("K132").Value = 111526.62
$ws_WVR.Range("M132").Value = -108996.62
